$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 252
$ws.Range("I28").Value = 276.44446
$ws.Range("K28").Value = 276.44446
$ws.Range("M28").Value = 208.55554
$ws.Range("H40").Value = 1914
$ws.Range("J40").Value = 2066.6667
$ws.Range("L40").Value = 2066.6667
$ws.Range("N40").Value = -2416.6667
$ws.Range("H70").Value = 34000
$ws.Range("I70").Value = 50250
$ws.Range("K70").Value = 150750
$ws.Range("M70").Value = -150480
$ws.Range("H73").Value = 34000
$ws.Range("I73").Value = 50250
$ws.Range("K73").Value = 150750
$ws.Range("M73").Value = -149814
$ws.Range("H101").Value = 1589.3334
$ws.Range("I101").Value = 978.6667
$ws.Range("J101").Value = 2200
$ws.Range("K101").Value = 2936.0001
$ws.Range("L101").Value = 6600
$ws.Range("M101").Value = -1314.0001
$ws.Range("N101").Value = -9844
$ws.Range("H137").Value = 19556.21
$ws.Range("I137").Value = 2257.1936
$ws.Range("J137").Value = 40181.96
$ws.Range("K137").Value = 6771.5808
$ws.Range("L137").Value = 120545.88
$ws.Range("M137").Value = -4221.5808
$ws.Range("N137").Value = -125645.88
$ws.Range("H138").Value = 2654.7144
$ws.Range("I138").Value = 1652.9445
$ws.Range("J138").Value = 3406.0417
$ws.Range("K138").Value = 4958.833500000001
$ws.Range("L138").Value = 10218.1251
$ws.Range("M138").Value = 181.1664999999994
$ws.Range("N138").Value = -20498.1251
$ws.Range("H141").Value = 2014.9286
$ws.Range("I141").Value = 1416.72
$ws.Range("K141").Value = 4250.16
$ws.Range("M141").Value = 929.8400000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17904.121
$ws.Range("I32").Value = 19397.508
$ws.Range("J32").Value = 8446
$ws.Range("K32").Value = 19397.508
$ws.Range("L32").Value = 8446
$ws.Range("M32").Value = -19110.508
$ws.Range("N32").Value = -9020
$ws.Range("H45").Value = 3540.842
$ws.Range("I45").Value = 2852
$ws.Range("J45").Value = 3858.7693
$ws.Range("K45").Value = 2852
$ws.Range("L45").Value = 3858.7693
$ws.Range("M45").Value = -2475
$ws.Range("N45").Value = -4612.7693
$ws.Range("H62").Value = 27500
$ws.Range("J62").Value = 27500
$ws.Range("L62").Value = 27500
$ws.Range("N62").Value = -28748
$ws.Range("H65").Value = 27500
$ws.Range("J65").Value = 27500
$ws.Range("L65").Value = 82500
$ws.Range("N65").Value = -88740
$ws.Range("H97").Value = 1792.2222
$ws.Range("I97").Value = 1930
$ws.Range("K97").Value = 1930
$ws.Range("M97").Value = -1434
$ws.Range("H102").Value = 5501.4287
$ws.Range("I102").Value = 3702
$ws.Range("K102").Value = 3702
$ws.Range("M102").Value = -2080
$ws.Range("H132").Value = 27767.35
$ws.Range("I132").Value = 2177.9167
$ws.Range("J132").Value = 66151.5
$ws.Range("K132").Value = 6533.750100000001
$ws.Range("L132").Value = 198454.5
$ws.Range("M132").Value = -4003.750100000001
$ws.Range("N132").Value = -203514.5
$ws.Range("H133").Value = 42726.816
$ws.Range("J133").Value = 42726.816
$ws.Range("L133").Value = 42726.816
$ws.Range("N133").Value = -47786.816

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3102.7144
$ws.Range("J94").Value = 4833.1665
$ws.Range("L94").Value = 4833.1665
$ws.Range("N94").Value = -5735.1665
$ws.Range("H134").Value = 56849.367
$ws.Range("I134").Value = 59840.223
$ws.Range("J134").Value = 3014
$ws.Range("K134").Value = 179520.669
$ws.Range("L134").Value = 9042
$ws.Range("M134").Value = -176985.669
$ws.Range("N134").Value = -14112

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12931.41
$ws.Range("I31").Value = 22076.895
$ws.Range("J31").Value = 4243.2
$ws.Range("K31").Value = 22076.895
$ws.Range("L31").Value = 4243.2
$ws.Range("M31").Value = -21781.895
$ws.Range("N31").Value = -4833.2
$ws.Range("H34").Value = 12931.41
$ws.Range("I34").Value = 22076.895
$ws.Range("J34").Value = 4243.2
$ws.Range("K34").Value = 22076.895
$ws.Range("L34").Value = 4243.2
$ws.Range("M34").Value = -21874.895
$ws.Range("N34").Value = -4647.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 841.86664
$ws.Range("I23").Value = 850
$ws.Range("J23").Value = 840.61536
$ws.Range("K23").Value = 2550
$ws.Range("L23").Value = 2521.84608
$ws.Range("M23").Value = -2315
$ws.Range("N23").Value = -2991.84608
$ws.Range("H35").Value = 1002
$ws.Range("I35").Value = 1002
$ws.Range("K35").Value = 3006
$ws.Range("M35").Value = -2718
$ws.Range("H38").Value = 179.8
$ws.Range("I38").Value = 233
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 699
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = -352
$ws.Range("N38").Value = -994
$ws.Range("H63").Value = 3679.875
$ws.Range("I63").Value = 1982.2
$ws.Range("J63").Value = 6509.3335
$ws.Range("K63").Value = 5946.6
$ws.Range("L63").Value = 19528.0005
$ws.Range("M63").Value = -5197.6
$ws.Range("N63").Value = -21026.0005
$ws.Range("H66").Value = 3679.875
$ws.Range("I66").Value = 1982.2
$ws.Range("J66").Value = 6509.3335
$ws.Range("K66").Value = 17839.8
$ws.Range("L66").Value = 58584.0015
$ws.Range("M66").Value = -14095.8
$ws.Range("N66").Value = -66072.0015
$ws.Range("H122").Value = 1141.875
$ws.Range("I122").Value = 387.875
$ws.Range("J122").Value = 1518.875
$ws.Range("K122").Value = 3490.875
$ws.Range("L122").Value = 13669.875
$ws.Range("M122").Value = -1040.875
$ws.Range("N122").Value = -18569.875
$ws.Range("H131").Value = 104973.49
$ws.Range("J131").Value = 114452.8
$ws.Range("L131").Value = 343358.4
$ws.Range("N131").Value = -353438.4
$ws.Range("H137").Value = 8447.700000000001
$ws.Range("I137").Value = 99999
$ws.Range("J137").Value = 3629.2104
$ws.Range("K137").Value = 299997
$ws.Range("L137").Value = 10887.6312
$ws.Range("M137").Value = -294897
$ws.Range("N137").Value = -21087.6312

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 25004500
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 33338334
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 33338334
$ws.Range("M58").Value = -2723
$ws.Range("N58").Value = -33338888
$ws.Range("H102").Value = 1862.4117
$ws.Range("I102").Value = 1784.0667
$ws.Range("J102").Value = 2450
$ws.Range("K102").Value = 1784.0667
$ws.Range("L102").Value = 2450
$ws.Range("M102").Value = -162.0667000000001
$ws.Range("N102").Value = -5694

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5344.2856
$ws.Range("I7").Value = 5238.421
$ws.Range("J7").Value = 6350
$ws.Range("K7").Value = 5238.421
$ws.Range("L7").Value = 6350
$ws.Range("M7").Value = -5126.421
$ws.Range("N7").Value = -6574
$ws.Range("H14").Value = 2766.6667
$ws.Range("J14").Value = 2766.6667
$ws.Range("L14").Value = 2766.6667
$ws.Range("N14").Value = -3110.6667
$ws.Range("H22").Value = 2128.5715
$ws.Range("J22").Value = 1449.5
$ws.Range("L22").Value = 1449.5
$ws.Range("N22").Value = -2039.5
$ws.Range("H27").Value = 2128.5715
$ws.Range("J27").Value = 1449.5
$ws.Range("L27").Value = 1449.5
$ws.Range("N27").Value = -1663.5
$ws.Range("H40").Value = 72337.125
$ws.Range("I40").Value = 125578.78
$ws.Range("K40").Value = 125578.78
$ws.Range("M40").Value = -125442.78
$ws.Range("H46").Value = 1885.5454
$ws.Range("I46").Value = 1737.8889
$ws.Range("J46").Value = 2550
$ws.Range("K46").Value = 1737.8889
$ws.Range("L46").Value = 2550
$ws.Range("M46").Value = -1549.8889
$ws.Range("N46").Value = -2926
$ws.Range("H126").Value = 5344.2856
$ws.Range("I126").Value = 5238.421
$ws.Range("J126").Value = 6350
$ws.Range("K126").Value = 15715.263
$ws.Range("L126").Value = 19050
$ws.Range("M126").Value = -13245.263
$ws.Range("N126").Value = -23990

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 39390
$ws.Range("J108").Value = 39390
$ws.Range("L108").Value = 39390
$ws.Range("N108").Value = -47070
